$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.2050359712230216
$ws.Cells.Item(2, 3).Value = 0.512589928057554
$ws.Cells.Item(2, 10).Value = 0.02158273381294964
$ws.Cells.Item(2, 16).Value = 0.1474820143884892
$ws.Cells.Item(2, 19).Value = 0.1133093525179856

$ws.Cells.Item(3, 2).Value = 0.0135593220338983
$ws.Cells.Item(3, 3).Value = 0.03050847457627119
$ws.Cells.Item(3, 10).Value = 0.03728813559322034
$ws.Cells.Item(3, 16).Value = 0.7389830508474576
$ws.Cells.Item(3, 19).Value = 0.1796610169491525

$ws.Cells.Item(4, 10).Value = 0.07608695652173914
$ws.Cells.Item(4, 15).Value = 0.0108695652173913
$ws.Cells.Item(4, 16).Value = 0.6413043478260869
$ws.Cells.Item(4, 19).Value = 0.2717391304347826

$ws.Cells.Item(6, 2).Value = 0.04646464646464647
$ws.Cells.Item(6, 4).Value = 0.01414141414141414
$ws.Cells.Item(6, 6).Value = 0.07474747474747474
$ws.Cells.Item(6, 10).Value = 0.2545454545454545
$ws.Cells.Item(6, 15).Value = 0.0303030303030303
$ws.Cells.Item(6, 17).Value = 0.1535353535353535
$ws.Cells.Item(6, 18).Value = 0.07474747474747474
$ws.Cells.Item(6, 19).Value = 0.3515151515151515

$ws.Cells.Item(7, 2).Value = 0.09975062344139651
$ws.Cells.Item(7, 4).Value = 0.02244389027431421
$ws.Cells.Item(7, 5).Value = 0.004987531172069825
$ws.Cells.Item(7, 6).Value = 0.06234413965087282
$ws.Cells.Item(7, 10).Value = 0.1371571072319202
$ws.Cells.Item(7, 15).Value = 0.02743142144638404
$ws.Cells.Item(7, 17).Value = 0.1745635910224439
$ws.Cells.Item(7, 18).Value = 0.08478802992518704
$ws.Cells.Item(7, 19).Value = 0.3865336658354115

$ws.Cells.Item(8, 2).Value = 0.08568075117370892
$ws.Cells.Item(8, 4).Value = 0.01643192488262911
$ws.Cells.Item(8, 5).Value = 0.001173708920187793
$ws.Cells.Item(8, 6).Value = 0.06455399061032864
$ws.Cells.Item(8, 10).Value = 0.1208920187793427
$ws.Cells.Item(8, 15).Value = 0.01995305164319249
$ws.Cells.Item(8, 17).Value = 0.1948356807511737
$ws.Cells.Item(8, 18).Value = 0.1009389671361502
$ws.Cells.Item(8, 19).Value = 0.3955399061032864

$ws.Cells.Item(9, 2).Value = 0.09213483146067415
$ws.Cells.Item(9, 4).Value = 0.008988764044943821
$ws.Cells.Item(9, 5).Value = 0.002247191011235955
$ws.Cells.Item(9, 6).Value = 0.0898876404494382
$ws.Cells.Item(9, 10).Value = 0.1280898876404494
$ws.Cells.Item(9, 15).Value = 0.0449438202247191
$ws.Cells.Item(9, 17).Value = 0.1685393258426966
$ws.Cells.Item(9, 18).Value = 0.1280898876404494
$ws.Cells.Item(9, 19).Value = 0.3370786516853932

$ws.Cells.Item(10, 2).Value = 0.09019044196909809
$ws.Cells.Item(10, 4).Value = 0.02048149478979518
$ws.Cells.Item(10, 5).Value = 0.0007186489399928135
$ws.Cells.Item(10, 6).Value = 0.06970894717930291
$ws.Cells.Item(10, 10).Value = 0.1189363995688106
$ws.Cells.Item(10, 15).Value = 0.02479338842975207
$ws.Cells.Item(10, 17).Value = 0.2091268415379087
$ws.Cells.Item(10, 18).Value = 0.1092346388789077
$ws.Cells.Item(10, 19).Value = 0.3568091987064319

$ws.Cells.Item(11, 7).Value = 0.1482649842271293
$ws.Cells.Item(11, 10).Value = 0.1025236593059937
$ws.Cells.Item(11, 11).Value = 0.194006309148265
$ws.Cells.Item(11, 12).Value = 0.5425867507886435
$ws.Cells.Item(11, 19).Value = 0.01261829652996845

$ws.Cells.Item(12, 7).Value = 0.7304347826086957
$ws.Cells.Item(12, 10).Value = 0.2492753623188406
$ws.Cells.Item(12, 11).Value = 0.005797101449275362
$ws.Cells.Item(12, 12).Value = 0.002898550724637681
$ws.Cells.Item(12, 19).Value = 0.01159420289855072

$ws.Cells.Item(13, 7).Value = 0.696969696969697
$ws.Cells.Item(13, 10).Value = 0.2525252525252525
$ws.Cells.Item(13, 19).Value = 0.0505050505050505

$ws.Cells.Item(15, 6).Value = 0.02240325865580448
$ws.Cells.Item(15, 8).Value = 0.1384928716904277
$ws.Cells.Item(15, 9).Value = 0.06313645621181263
$ws.Cells.Item(15, 10).Value = 0.3258655804480652
$ws.Cells.Item(15, 11).Value = 0.06720977596741344
$ws.Cells.Item(15, 13).Value = 0.01832993890020367
$ws.Cells.Item(15, 15).Value = 0.07535641547861507
$ws.Cells.Item(15, 19).Value = 0.2892057026476578

$ws.Cells.Item(16, 6).Value = 0.03107344632768362
$ws.Cells.Item(16, 8).Value = 0.1581920903954802
$ws.Cells.Item(16, 9).Value = 0.07344632768361582
$ws.Cells.Item(16, 10).Value = 0.4209039548022599
$ws.Cells.Item(16, 11).Value = 0.1271186440677966
$ws.Cells.Item(16, 13).Value = 0.02542372881355932
$ws.Cells.Item(16, 15).Value = 0.03389830508474576
$ws.Cells.Item(16, 19).Value = 0.1299435028248588

$ws.Cells.Item(17, 6).Value = 0.01776384535005225
$ws.Cells.Item(17, 8).Value = 0.1630094043887147
$ws.Cells.Item(17, 9).Value = 0.1086729362591431
$ws.Cells.Item(17, 10).Value = 0.4263322884012539
$ws.Cells.Item(17, 11).Value = 0.08986415882967608
$ws.Cells.Item(17, 13).Value = 0.02716823406478579
$ws.Cells.Item(17, 14).Value = 0.001044932079414838
$ws.Cells.Item(17, 15).Value = 0.06374085684430512
$ws.Cells.Item(17, 19).Value = 0.1024033437826541

$ws.Cells.Item(18, 6).Value = 0.02707930367504836
$ws.Cells.Item(18, 8).Value = 0.1972920696324952
$ws.Cells.Item(18, 9).Value = 0.1005802707930367
$ws.Cells.Item(18, 10).Value = 0.4119922630560928
$ws.Cells.Item(18, 11).Value = 0.08123791102514506
$ws.Cells.Item(18, 13).Value = 0.01740812379110251
$ws.Cells.Item(18, 15).Value = 0.05029013539651837
$ws.Cells.Item(18, 19).Value = 0.1141199226305609

$ws.Cells.Item(19, 6).Value = 0.01427469135802469
$ws.Cells.Item(19, 8).Value = 0.1809413580246914
$ws.Cells.Item(19, 9).Value = 0.09027777777777778
$ws.Cells.Item(19, 10).Value = 0.3981481481481481
$ws.Cells.Item(19, 11).Value = 0.1091820987654321
$ws.Cells.Item(19, 13).Value = 0.02006172839506173
$ws.Cells.Item(19, 14).Value = 0.0003858024691358024
$ws.Cells.Item(19, 15).Value = 0.06404320987654322
$ws.Cells.Item(19, 19).Value = 0.1226851851851852
